# feat: add 2022-Q3 data
#
# Before:  总计 (sheet1) | 2022-Q2 (sheet2, fund holdings table)
# After:   总计 (sheet1) | 2022-Q3 (sheet2, NEW fund holdings table) | 2022-Q2 (sheet3, old fund holdings table)
#
# Excel always assigns a brand-new sheet the next free sheetId, so to land on
# sheetId 2 / rId2 for "2022-Q3" (matching the target file) we repurpose the
# existing "2022-Q2" sheet in place (rename + replace its data) and create a
# *new* sheet for "2022-Q2"'s old content, appended right after it.

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item("总计")
$oldQ2   = $wb.Worksheets.Item("2022-Q2")

# ---------------------------------------------------------------------------
# 1) "总计": insert a new row 2 for the 2022-Q3 summary, push 2022-Q2 to row 3
# ---------------------------------------------------------------------------
$summary.Rows.Item(2).Insert()

# The freshly inserted row inherits stray formatting from the row above it;
# reset it back to the default (unstyled) look using an untouched cell.
$summary.Range("F10").Copy()
$summary.Range("A2:D2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 0.09

# Give the new index cell (A2) the same bold/border/center-top look as the
# rest of the "总计" sheet's label cells.
$summary.Range("B1").Copy()
$summary.Range("A2").PasteSpecial(-4122)

# The row that got pushed down keeps its old index value (0); it should now
# read 1, since it is the second entry in the list.
$summary.Range("A3").Value = 1

# ---------------------------------------------------------------------------
# 2) Turn the current "2022-Q2" sheet into "2022-Q3" and give it the new data
# ---------------------------------------------------------------------------
$oldQ2.Name = "2022-Q3"

# Clear out the previous 5-row fund table before writing the 3-row one.
$oldQ2.Range("A1:H5").Clear()

$oldQ2.Range("B1").Value = "基金代码"
$oldQ2.Range("C1").Value = "基金名称"
$oldQ2.Range("D1").Value = "基金规模"
$oldQ2.Range("E1").Value = "股票总仓位"
$oldQ2.Range("F1").Value = "仓位占比"
$oldQ2.Range("G1").Value = "持有市值(亿元)"
$oldQ2.Range("H1").Value = "仓位排名"
$oldQ2.Range("B1:H1").Font.Bold = $true
$oldQ2.Range("B1:H1").Borders.LineStyle = 1
$oldQ2.Range("B1:H1").HorizontalAlignment = -4108
$oldQ2.Range("B1:H1").VerticalAlignment = -4160

$oldQ2.Range("A2").Value = 0
$oldQ2.Range("H2").Value = 8

$oldQ2.Range("A3").Value = 1
$oldQ2.Range("H3").Value = 8

$oldQ2.Range("A2:A3").Font.Bold = $true
$oldQ2.Range("A2:A3").Borders.LineStyle = 1
$oldQ2.Range("A2:A3").HorizontalAlignment = -4108
$oldQ2.Range("A2:A3").VerticalAlignment = -4160

# B/C/D/E/F/G on the data rows are stored as text, not numbers (fund codes
# like "013340" need their leading zero kept) - force text type via a "@"
# number format, then immediately paste a blank cell's format back over the
# range so no stray numFmt style is left behind on the cells.
$oldQ2.Range("B2:G3").NumberFormat = "@"
$oldQ2.Range("B2").Value = "013340"
$oldQ2.Range("C2").Value = "创金合信芯片产业股票C"
$oldQ2.Range("D2").Value = "0.93"
$oldQ2.Range("E2").Value = "90.74"
$oldQ2.Range("F2").Value = "4.72"
$oldQ2.Range("G2").Value = "0.0439"
$oldQ2.Range("B3").Value = "013339"
$oldQ2.Range("C3").Value = "创金合信芯片产业股票A"
$oldQ2.Range("D3").Value = "0.91"
$oldQ2.Range("E3").Value = "90.74"
$oldQ2.Range("F3").Value = "4.72"
$oldQ2.Range("G3").Value = "0.0430"

$oldQ2.Range("J1").Copy()
$oldQ2.Range("B2:G3").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3) Add a brand-new sheet, named "2022-Q2", right after "2022-Q3" and give
#    it the fund table that used to live in the "2022-Q2" sheet.
# ---------------------------------------------------------------------------
$newQ2 = $wb.Worksheets.Add($null, $oldQ2)
$newQ2.Name = "2022-Q2"

$newQ2.Range("B1").Value = "基金代码"
$newQ2.Range("C1").Value = "基金名称"
$newQ2.Range("D1").Value = "基金规模"
$newQ2.Range("E1").Value = "股票总仓位"
$newQ2.Range("F1").Value = "仓位占比"
$newQ2.Range("G1").Value = "持有市值(亿元)"
$newQ2.Range("H1").Value = "仓位排名"
$newQ2.Range("B1:H1").Font.Bold = $true
$newQ2.Range("B1:H1").Borders.LineStyle = 1
$newQ2.Range("B1:H1").HorizontalAlignment = -4108
$newQ2.Range("B1:H1").VerticalAlignment = -4160

$newQ2.Range("A2").Value = 0
$newQ2.Range("H2").Value = 8

$newQ2.Range("A3").Value = 1
$newQ2.Range("H3").Value = 10

$newQ2.Range("A4").Value = 2
$newQ2.Range("H4").Value = 10

$newQ2.Range("A5").Value = 3
$newQ2.Range("H5").Value = 8

$newQ2.Range("A2:A5").Font.Bold = $true
$newQ2.Range("A2:A5").Borders.LineStyle = 1
$newQ2.Range("A2:A5").HorizontalAlignment = -4108
$newQ2.Range("A2:A5").VerticalAlignment = -4160

# B/C/D/E/F/G on the data rows are stored as text, not numbers (fund codes
# like "013273"/"007804" need their leading zero kept) - force text type via
# a "@" number format, then paste a blank cell's format back over the range
# so no stray numFmt style is left behind on the cells.
$newQ2.Range("B2:G5").NumberFormat = "@"
$newQ2.Range("B2").Value = "310318"
$newQ2.Range("C2").Value = "申万菱信沪深300指数增强A"
$newQ2.Range("D2").Value = "9.56"
$newQ2.Range("E2").Value = "88.14"
$newQ2.Range("F2").Value = "0.03"
$newQ2.Range("G2").Value = "0.0029"
$newQ2.Range("B3").Value = "013273"
$newQ2.Range("C3").Value = "招商沪深300地产等权重指数C"
$newQ2.Range("D3").Value = "12.33"
$newQ2.Range("E3").Value = "94.91"
$newQ2.Range("F3").Value = "0.01"
$newQ2.Range("G3").Value = "0.0012"
$newQ2.Range("B4").Value = "161721"
$newQ2.Range("C4").Value = "招商沪深300地产等权重指数"
$newQ2.Range("D4").Value = "5.87"
$newQ2.Range("E4").Value = "94.91"
$newQ2.Range("F4").Value = "0.01"
$newQ2.Range("G4").Value = "0.0006"
$newQ2.Range("B5").Value = "007804"
$newQ2.Range("C5").Value = "申万菱信沪深300指数增强C"
$newQ2.Range("D5").Value = "1.18"
$newQ2.Range("E5").Value = "88.14"
$newQ2.Range("F5").Value = "0.03"
$newQ2.Range("G5").Value = "0.0004"

$newQ2.Range("J1").Copy()
$newQ2.Range("B2:G5").PasteSpecial(-4122)

$summary.Select()
$summary.Range("A1").Select()
